$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# The document opened with two back-to-back empty paragraphs at the very
# top: a plain one, followed by one carrying <w:vanish/><w:specVanish/>
# run formatting. Remove the redundant leading (non-vanished) paragraph so
# the vanished one becomes the document's first paragraph.
$firstPara = $d.Paragraphs.First
$firstPara.Range.Delete()

# --- Change 2 ---------------------------------------------------------
# Word stamps a "_GoBack" bookmark at the location of the most recent
# edit/selection. Relocate it from the last (empty) paragraph at the end
# of the document to the start of the "12. ..." milestone paragraph.
# Adding a bookmark named "_GoBack" again simply relocates the existing
# one instead of creating a duplicate, mirroring how Word itself tracks
# this.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("12. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $target = $searchRange.Paragraphs.First
    $bmStart = $target.Range.Start
    $bmRange = $d.Range($bmStart, $bmStart)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
